# Weekly update: insert a new price-record row for "Zapallo italiano" at
# Terminal La Palmera de La Serena, pushing the existing rows 134-212 down
# to 135-213.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(134).Insert()

$ws.Range("A134").Value = 8
$ws.Range("B134").Value = "Terminal La Palmera de La Serena"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44518
$ws.Range("E134").Value = 4
$ws.Range("F134").Value = 100112032
$ws.Range("G134").Value = "Zapallo italiano"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 400
$ws.Range("K134").Value = 10000
$ws.Range("L134").Value = 11000
$ws.Range("M134").Value = 10500
$ws.Range("N134").Value = "`$/caja 70 unidades"
$ws.Range("O134").Value = "Provincia de Limarí"
$ws.Range("P134").Value = 150
$ws.Range("Q134").Value = 70
$ws.Range("R134").Value = "Hortaliza"
